$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet contained two "POC" rows that injected extra shared-string
# values ("ABC" / "DEF") into J1, J2, B9 and B10. Remove that POC data.
$ws.Range("J1").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()

# Leave the UI selection where it ended up after the cleanup.
$ws.Range("J2").Select()
